$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: clear description and amount (B3, D3) to simulate invalid/missing data
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()

# New rows 6-9 with additional (invalid) test data
$ws.Range("A6").Value = "A1004"
$ws.Range("C6").Value = "GBP"
$ws.Range("D6").Value = 10000

$ws.Range("A7").Value = "A1005"
$ws.Range("B7").Value = "Description for A1001"
$ws.Range("C7").Value = "USD"

$ws.Range("A8").Value = "A1006"
$ws.Range("C8").Value = "INR"
$ws.Range("D8").Value = 10000

$ws.Range("A9").Value = "A1007"
$ws.Range("B9").Value = "Description for A1003"
$ws.Range("C9").Value = "AUD"
$ws.Range("D9").Value = 10000

# Update the selection to match the final state
$ws.Range("D7").Select()
